$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artifacts")

# ---------------------------------------------------------------------
# 1. Make room for the new MODS columns that describe the source book
#    the new Tiyo Soga image was published in. Two pairs of columns are
#    inserted (matching the shift AA->AC, AB->AD, AC->AG, AD->AH, ...
#    seen between the old and new workbook layouts).
# ---------------------------------------------------------------------
$ws.Range("AA1:AB1").EntireColumn.Insert()
$ws.Range("AE1:AF1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. Append the new artifact row (row 74) describing the "Tiyo Soga
#    (with Facsimile Signature)" image.
# ---------------------------------------------------------------------
$ws.Range("A74").Value = "https://archive.org/details/tiyosogapageofso00chal/page/n7/mode/2up"
$ws.Range("B74").Value = "Public domain"
$ws.Range("C74").Value = "Tiyo Soga seat in quarter profile facing forward and holding a book."
$ws.Range("F74").Value = "liv_021070"
$ws.Range("G74").Value = "Tiyo Soga (with Facsimile Signature)"
$ws.Range("H74").Value = "Tiyo Soga (with Facsimile Signature), [late nineteenth century]"

# I74 carries the same "not-left-aligned" style used by other Anonymous
# name cells further up the sheet (e.g. I71) - copy its format first.
$ws.Range("I71").Copy()
$ws.Range("I74").PasteSpecial(-4104)
$ws.Range("I74").Value = "Anonymous"

$ws.Range("L74").Value = "photographs"
$ws.Range("M74").Value = " publications (documents)"

# V74 uses the "author" style (style 10) used by other author cells
# (e.g. V62) - copy its format first.
$ws.Range("V62").Copy()
$ws.Range("V74").PasteSpecial(-4104)
$ws.Range("V74").Value = "Chalmers, John A."

$ws.Range("X74").Value = "Tiyo Soga: A Page of South African Mission Work"
$ws.Range("Y74").Value = "Edinburgh"
$ws.Range("Z74").Value = "London"
$ws.Range("AA74").Value = "Glasgow"
$ws.Range("AB74").Value = "Grahamstown, Cape Colony"
$ws.Range("AC74").Value = "Andrew Elliot"
$ws.Range("AD74").Value = "Hodder & Stoughton"
$ws.Range("AE74").Value = "David Bryce & Son"
$ws.Range("AF74").Value = "James Kay"
$ws.Range("AG74").Value = 1877
$ws.Range("AH74").Value = "frontispiece"

$ws.Rows.Item(74).RowHeight = 68

# ---------------------------------------------------------------------
# 3. Scroll/select the Artifacts sheet the way the author left it and
#    make it the active tab (previously "Documents" was active).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A70").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("AS74").Select()
$excel.ActiveWindow.ScrollColumn = 34
